$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.1052853234940633
$ws.Range("C2").Value = 0.4900476106607171
$ws.Range("D2").Value = 0.6183260505831405
$ws.Range("E2").Value = 0.7863371100127098
$ws.Range("F2").Value = 0.7868594147108519
$ws.Range("G2").Value = 52

$ws.Range("B3").Value = 0.3646584415625562
$ws.Range("C3").Value = 0.722637768445621
$ws.Range("D3").Value = 1.426035144228675
$ws.Range("E3").Value = 1.194167134126825
$ws.Range("F3").Value = 1.148442664015213
$ws.Range("G3").Value = 51

$ws.Range("B4").Value = 0.4532797608695453
$ws.Range("C4").Value = 0.9902620268281389
$ws.Range("D4").Value = 3.506287014911406
$ws.Range("E4").Value = 1.872508214911595
$ws.Range("F4").Value = 1.835262443918751
$ws.Range("G4").Value = 50

$ws.Range("B5").Value = 0.4187269926694491
$ws.Range("C5").Value = 1.116030423578108
$ws.Range("D5").Value = 4.462820609660299
$ws.Range("E5").Value = 2.112538901336564
$ws.Range("F5").Value = 2.092082930599334
$ws.Range("G5").Value = 49

$ws.Range("B6").Value = 0.2603619109577968
$ws.Range("C6").Value = 0.935772708153241
$ws.Range("D6").Value = 3.863674849690833
$ws.Range("E6").Value = 1.965623272575606
$ws.Range("F6").Value = 1.96892104161102
$ws.Range("G6").Value = 48

$ws.Range("B7").Value = 0.2925042579902095
$ws.Range("C7").Value = 0.9789877322222572
$ws.Range("D7").Value = 4.769865589305868
$ws.Range("E7").Value = 2.184002195352804
$ws.Range("F7").Value = 2.19502714946493
$ws.Range("G7").Value = 36

$ws.Range("B8").Value = 0.2341597008279832
$ws.Range("C8").Value = 0.9857524095050276
$ws.Range("D8").Value = 4.861337304003894
$ws.Range("E8").Value = 2.204844054350306
$ws.Range("F8").Value = 2.224381796787603
$ws.Range("G8").Value = 35

$ws.Range("B9").Value = 0.1803491558900733
$ws.Range("C9").Value = 1.423780144984279
$ws.Range("D9").Value = 8.982173702935924
$ws.Range("E9").Value = 2.997027477841324
$ws.Range("F9").Value = 3.078327104205846
$ws.Range("G9").Value = 18

$ws.Range("B10").Value = -0.6058116565806465
$ws.Range("C10").Value = 1.074425645603604
$ws.Range("D10").Value = 6.535069776881333
$ws.Range("E10").Value = 2.556378253874284
$ws.Range("F10").Value = 2.60477795886628
$ws.Range("G10").Value = 11
